$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "ODI Bowling Extra"
$ws.Range("A1").Value = "MATCH_CODE"
